$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.304.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.217.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.62"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.23%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.548.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.216.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.244.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.65%  "
$ws.Range("E22").Value = "  -7.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.95%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.34%  "
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.25%  "
$ws.Range("E29").Value = "  -4.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.43%  "
$ws.Range("E33").Value = "  -7.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0718"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.57%  "
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.56%  "
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0278"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("E40").Value = "  -5.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.200"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.57%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.28%  "
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("E50").Value = "  -5.57%  "
$ws.Range("E51").Value = "  -5.39%  "
